$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: update event name, attendee count, and min ticket price (now "not for sale")
    $ws.Range("C2").Value = "南宁·熊喵M动漫嘉年华·万圣派对（取消）"
    $ws.Range("F2").Value = 336
    $ws.Range("G2").Value = "不可售"

    # Row 3: update attendee count only
    $ws.Range("F3").Value = 1372
}
